$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range('A32').Value = 'Indonesia'
$ws.Range('A33').Value = 'Emiratos Arabes Unidos'
$ws.Range('A39').Value = 'Polonia'
$ws.Range('A40').Value = 'Suiza'
$ws.Range('A144').Value = 'Estado de Palestina'
$ws.Range('A145').Value = 'Jamaica'
$ws.Range('A202').Value = 'Dominica'
$ws.Range('A203').Value = 'Fiyi'
$ws.Range('A206').Value = 'Islas Malvinas'
$ws.Range('A207').Value = 'Groenlandia'
$ws.Range('A210').Value = 'Seychelles'
$ws.Range('A211').Value = 'Montserrat'
$ws.Range('A213').Value = 'Papua Nueva Guinea'
$ws.Range('A214').Value = 'Islas Virgenes Britanicas'
$ws.Range('A107').Value = 'Albania'
$ws.Range('A108').Value = 'Sudan del Sur'
$ws.Range('A109').Value = 'Nicaragua'
$ws.Range('A110').Value = 'Islandia'
$ws.Range('A111').Value = 'Lituania'
$ws.Range('B4').Value = 2263756
$ws.Range('C4').Value = 105
$ws.Range('D4').Value = 931079
$ws.Range('E4').Value = 1211989
$ws.Range('B7').Value = 381539
$ws.Range('C7').Value = 448
$ws.Range('D7').Value = 205245
$ws.Range('E7').Value = 163688
$ws.Range('G7').Value = 2
$ws.Range('H7').Value = 12606
$ws.Range('B20').Value = 105535
$ws.Range('C20').Value = 3243
$ws.Range('D20').Value = 42945
$ws.Range('E20').Value = 61202
$ws.Range('G20').Value = 45
$ws.Range('H20').Value = 1388
$ws.Range('B25').Value = 60476
$ws.Range('C25').Value = 128
$ws.Range('D25').Value = 16751
$ws.Range('E25').Value = 34030
$ws.Range('G25').Value = 12
$ws.Range('H25').Value = 9695
$ws.Range('B32').Value = 43803
$ws.Range('C32').Value = 1041
$ws.Range('D32').Value = 17349
$ws.Range('E32').Value = 24081
$ws.Range('G32').Value = 34
$ws.Range('H32').Value = 2373
$ws.Range('B33').Value = 43752
$ws.Range('D33').Value = 30241
$ws.Range('E33').Value = 13213
$ws.Range('H33').Value = 298
$ws.Range('B39').Value = 31316
$ws.Range('C39').Value = 301
$ws.Range('D39').Value = 15698
$ws.Range('E39').Value = 14284
$ws.Range('G39').Value = 18
$ws.Range('H39').Value = 1334
$ws.Range('B40').Value = 31200
$ws.Range('D40').Value = 28900
$ws.Range('E40').Value = 344
$ws.Range('H40').Value = 1956
$ws.Range('B42').Value = 27878
$ws.Range('C42').Value = 346
$ws.Range('D42').Value = 7962
$ws.Range('E42').Value = 19368
$ws.Range('G42').Value = 2
$ws.Range('H42').Value = 548
$ws.Range('B43').Value = 27670
$ws.Range('C43').Value = 852
$ws.Range('D43').Value = 13974
$ws.Range('E43').Value = 13571
$ws.Range('G43').Value = 6
$ws.Range('H43').Value = 125
$ws.Range('E50').Value = 5678
$ws.Range('G50').Value = 1
$ws.Range('H50').Value = 56
$ws.Range('B55').Value = 17271
$ws.Range('C55').Value = 48
$ws.Range('D55').Value = 16141
$ws.Range('E55').Value = 442
$ws.Range('D57').Value = 7525
$ws.Range('E57').Value = 5132
$ws.Range('G57').Value = 5
$ws.Range('H57').Value = 449
$ws.Range('B68').Value = 9280
$ws.Range('C68').Value = 206
$ws.Range('D68').Value = 8081
$ws.Range('E68').Value = 986
$ws.Range('B70').Value = 8535
$ws.Range('C70').Value = 6
$ws.Range('D70').Value = 8070
$ws.Range('E70').Value = 344
$ws.Range('D102').Value = 1759
$ws.Range('E102').Value = 370
$ws.Range('B107').Value = 1838
$ws.Range('C107').Value = 50
$ws.Range('D107').Value = 1114
$ws.Range('E107').Value = 682
$ws.Range('G107').Value = 3
$ws.Range('H107').Value = 42
$ws.Range('B108').Value = 1830
$ws.Range('D108').Value = 117
$ws.Range('E108').Value = 1681
$ws.Range('H108').Value = 32
$ws.Range('B109').Value = 1823
$ws.Range('D109').Value = 1238
$ws.Range('E109').Value = 521
$ws.Range('H109').Value = 64
$ws.Range('B110').Value = 1816
$ws.Range('C110').Value = 0
$ws.Range('D110').Value = 1801
$ws.Range('E110').Value = 5
$ws.Range('H110').Value = 10
$ws.Range('B111').Value = 1792
$ws.Range('C111').Value = 8
$ws.Range('D111').Value = 1462
$ws.Range('E111').Value = 254
$ws.Range('H111').Value = 76
$ws.Range('B115').Value = 1513
$ws.Range('C115').Value = 2
$ws.Range('E115').Value = 45
$ws.Range('B144').Value = 639
$ws.Range('C144').Value = 39
$ws.Range('D144').Value = 415
$ws.Range('E144').Value = 221
$ws.Range('H144').Value = 3
$ws.Range('B145').Value = 638
$ws.Range('C145').Value = 12
$ws.Range('D145').Value = 458
$ws.Range('E145').Value = 170
$ws.Range('H145').Value = 10
$ws.Range('D210').Value = 11
$ws.Range('H210').Value = 0
$ws.Range('D211').Value = 10
$ws.Range('H211').Value = 1
$ws.Range('D213').Value = 8
$ws.Range('H213').Value = 0
$ws.Range('D214').Value = 7
$ws.Range('H214').Value = 1
$ws.Range('A1').Value = 'Datos actualizados a 19 de Junio de 2020 a las 11:53'
